$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Q0) updated values
$ws.Range("B3").Value = 0.0243169190249313
$ws.Range("C3").Value = 0.5668789673031478
$ws.Range("D3").Value = 0.5577053810991515
$ws.Range("E3").Value = 0.7467967468455868
$ws.Range("F3").Value = 0.7668538649144974
$ws.Range("G3").Value = 19

# Row 4 (Q1) updated values
$ws.Range("B4").Value = 0.3353634677292027
$ws.Range("C4").Value = 0.6039326226106767
$ws.Range("D4").Value = 0.7960447219234688
$ws.Range("E4").Value = 0.8922133836271841
$ws.Range("F4").Value = 0.8507563830411711
$ws.Range("G4").Value = 18
